$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new device row (row 12): Samsung TV
$ws.Cells.Item(12, 2).Value = "Samsung"
$ws.Cells.Item(12, 3).Value = "TV"
$ws.Cells.Item(12, 4).Value = "10.1.77.107"

# Move the active selection down, as Excel does after entering data in a cell
[void]$ws.Range("D13").Select()
